$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2" = 1.02
    "C2" = 1.008199406363829
    "D2" = 1.028405560243696
    "E2" = 1.010736026748353
    "F2" = 1.022458468557891
    "I2" = 1.030470695751452
    "J2" = 1.013466621627579
    "K2" = 1.031222471817723
    "L2" = 1.013605317779909
    "M2" = 1.025292789681212
    "N2" = 1.00846097436131
    "B3" = 1.02
    "C3" = 1.009980992467433
    "D3" = 1.028888867513922
    "E3" = 1.012272231573513
    "F3" = 1.024316379258391
    "I3" = 1.030615051617936
    "J3" = 1.01487556979671
    "K3" = 1.031514583961475
    "L3" = 1.014943613307908
    "M3" = 1.026954501674668
    "N3" = 1.008944619702715
    "B4" = 1.02
    "C4" = 1.011130102647602
    "D4" = 1.029198749604333
    "E4" = 1.013263261496633
    "F4" = 1.025511196988755
    "I4" = 1.030704755590867
    "J4" = 1.015783417638934
    "K4" = 1.031700036525147
    "L4" = 1.015806118417662
    "M4" = 1.028021979428412
    "N4" = 1.009255767867436
    "B5" = 1.02
    "C5" = 1.011612322521877
    "D5" = 1.02932834109454
    "E5" = 1.013679187650658
    "F5" = 1.026011750915302
    "I5" = 1.030741580538863
    "J5" = 1.016164173168042
    "K5" = 1.031777147365479
    "L5" = 1.016167900484134
    "M5" = 1.028468905192911
    "N5" = 1.009386148030789
    "B6" = 1.02
    "C6" = 1.011693239110816
    "D6" = 1.029350060007496
    "E6" = 1.013748982692664
    "F6" = 1.026095694169603
    "I6" = 1.030747711620971
    "J6" = 1.016228051123685
    "K6" = 1.031790044559374
    "L6" = 1.016228597880755
    "M6" = 1.028543838448504
    "N6" = 1.009408014567595
    "B7" = 1.02
    "C7" = 1.011136549470888
    "D7" = 1.029200483896222
    "E7" = 1.013268821863812
    "F7" = 1.025517892256369
    "I7" = 1.030705251131012
    "J7" = 1.01578850884401
    "K7" = 1.031701070237484
    "L7" = 1.015810955751498
    "M7" = 1.028027958487598
    "N7" = 1.009257511682621
    "B8" = 1.02
    "C8" = 1.008802278088652
    "D8" = 1.028569485964282
    "E8" = 1.011255823903459
    "F8" = 1.023087896141609
    "I8" = 1.030520248513936
    "J8" = 1.013943585339947
    "K8" = 1.031321930184477
    "L8" = 1.014058326019438
    "M8" = 1.025855989795965
    "N8" = 1.008624800890921
    "B9" = 1.02
    "C9" = 1.004659852234417
    "D9" = 1.027435781528671
    "E9" = 1.007685041029768
    "F9" = 1.018748626990426
    "I9" = 1.030165872476062
    "J9" = 1.010662554869645
    "K9" = 1.030626552345071
    "L9" = 1.010942839536329
    "M9" = 1.021968527401275
    "N9" = 1.007495848148559
    "B10" = 1.02
    "C10" = 1.001877456530309
    "D10" = 1.026665337858934
    "E10" = 1.005287694533441
    "F10" = 1.015816080579863
    "I10" = 1.029910520004136
    "J10" = 1.008454070290377
    "K10" = 1.030144623004911
    "L10" = 1.008846737143951
    "M10" = 1.019335378657959
    "N10" = 1.006733452815206
    "B11" = 1.02
    "C11" = 1.000667453607731
    "D11" = 1.026328255826409
    "E11" = 1.004245414529502
    "F11" = 1.014536573055618
    "I11" = 1.029795411600501
    "J11" = 1.007492549719182
    "K11" = 1.029931588405758
    "L11" = 1.007934372454489
    "M11" = 1.018185112481438
    "N11" = 1.006400939250378
    "B12" = 1.02
    "C12" = 1.000217199172281
    "D12" = 1.026202526307366
    "E12" = 1.003857613917235
    "F12" = 1.01405982696987
    "I12" = 1.029751972418833
    "J12" = 1.00713459377455
    "K12" = 1.029851802908561
    "L12" = 1.007594750594229
    "M12" = 1.017756314669423
    "N12" = 1.006277063263189
    "B13" = 1.02
    "C13" = 1.000313817101804
    "D13" = 1.026229519352788
    "E13" = 1.003940828241598
    "F13" = 1.014162157977488
    "I13" = 1.029761321183172
    "J13" = 1.007211413237316
    "K13" = 1.029868946816384
    "L13" = 1.007667633901541
    "M13" = 1.017848363183346
    "N13" = 1.006303651732442
    "B14" = 1.02
    "C14" = 1.000630252002251
    "D14" = 1.026317873640243
    "E14" = 1.004213372206165
    "F14" = 1.014497195404279
    "I14" = 1.029791834838225
    "J14" = 1.007462977485511
    "K14" = 1.02992500667608
    "L14" = 1.007906314176563
    "M14" = 1.018149699413077
    "N14" = 1.00639070711888
    "B15" = 1.02
    "C15" = 1.00082511055636
    "D15" = 1.026372242422932
    "E15" = 1.004381208663192
    "F15" = 1.014703426196422
    "I15" = 1.029810544804343
    "J15" = 1.00761786737869
    "K15" = 1.029959460199703
    "L15" = 1.008053275857377
    "M15" = 1.018335158208146
    "N15" = 1.00644429617195
    "B16" = 1.02
    "C16" = 1.001957648761103
    "D16" = 1.026687635658844
    "E16" = 1.005356776891426
    "F16" = 1.015900791021025
    "I16" = 1.029918063665843
    "J16" = 1.008517771451075
    "K16" = 1.030158669591454
    "L16" = 1.008907186496333
    "M16" = 1.019411503552831
    "N16" = 1.006755469720891
    "B17" = 1.02
    "C17" = 1.002666649701268
    "D17" = 1.02688454309132
    "E17" = 1.005967584161648
    "F17" = 1.016649254201919
    "I17" = 1.029984291590269
    "J17" = 1.009080843856665
    "K17" = 1.030282461816951
    "L17" = 1.009441541466629
    "M17" = 1.020083949019422
    "N17" = 1.00695001596823
    "B18" = 1.02
    "C18" = 1.003079697114114
    "D18" = 1.026999060614117
    "E18" = 1.006323452459918
    "F18" = 1.017084886232811
    "I18" = 1.030022483323912
    "J18" = 1.00940877088749
    "K18" = 1.030354247477751
    "L18" = 1.009752765649112
    "M18" = 1.020475202467008
    "N18" = 1.007063261207401
    "B19" = 1.02
    "C19" = 1.003220451351786
    "D19" = 1.027038051269005
    "E19" = 1.006444726169426
    "F19" = 1.017233267795256
    "I19" = 1.030035431467616
    "J19" = 1.009520500665074
    "K19" = 1.030378653264921
    "L19" = 1.0098588082146
    "M19" = 1.020608445438216
    "N19" = 1.007101836036499
    "B20" = 1.02
    "C20" = 1.002590632581257
    "D20" = 1.026863451465139
    "E20" = 1.005902092371679
    "F20" = 1.016569047949643
    "I20" = 1.029977231260185
    "J20" = 1.009020483735633
    "K20" = 1.030269223554725
    "L20" = 1.009384257527853
    "M20" = 1.020011902738245
    "N20" = 1.006929166863467
    "B21" = 1.02
    "C21" = 1.000537092233659
    "D21" = 1.026291869923291
    "E21" = 1.004133132933141
    "F21" = 1.014398576278381
    "I21" = 1.029782868183062
    "J21" = 1.007388920438185
    "K21" = 1.029908516529494
    "L21" = 1.007836049046742
    "M21" = 1.018061006023993
    "N21" = 1.006365081622595
    "B22" = 1.02
    "C22" = 0.999241277965408
    "D22" = 1.025929472208775
    "E22" = 1.003017139724131
    "F22" = 1.013025339744179
    "I22" = 1.029656713477161
    "J22" = 1.006358428586246
    "K22" = 1.029677935966401
    "L22" = 1.006858402552158
    "M22" = 1.016825492230817
    "N22" = 1.006008300110387
    "B23" = 1.019999999999999
    "C23" = 0.99992866429198
    "D23" = 1.026121872608434
    "E23" = 1.00360911319645
    "F23" = 1.013754139128842
    "I23" = 1.029723965280537
    "J23" = 1.006905160026701
    "K23" = 1.029800530509603
    "L23" = 1.007377077839923
    "M23" = 1.017481312933195
    "N23" = 1.006197639724978
    "B24" = 1.02
    "C24" = 1.002624982992661
    "D24" = 1.026872982899426
    "E24" = 1.005931686544309
    "F24" = 1.016605292587852
    "I24" = 1.029980422872689
    "J24" = 1.009047759431762
    "K24" = 1.030275206654422
    "L24" = 1.009410143082756
    "M24" = 1.020044460356534
    "N24" = 1.006938588387941
    "B25" = 1.02
    "C25" = 1.005734343001851
    "D25" = 1.027731452884222
    "E25" = 1.00861106955114
    "F25" = 1.019877342222893
    "I25" = 1.030260852055462
    "J25" = 1.011514435009111
    "K25" = 1.030809557105707
    "L25" = 1.011751571290116
    "M25" = 1.022980767330699
    "N25" = 1.007789405759991
}

foreach ($key in $values.Keys) {
    $ws.Range($key).Value = $values[$key]
}

Write-Host "Updated $($values.Keys.Count) cells"